$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# Row 126: Jamshedpur FC vs FC Goa - match result now known (2-3, Away win)
# ---------------------------------------------------------------
$ws.Cells.Item(126, 8).Value = 2
$ws.Cells.Item(126, 9).Value = 3
$ws.Cells.Item(126, 10).Value = 'A'
$ws.Cells.Item(126, 14).Value = 5.5
$ws.Cells.Item(126, 15).Value = 4.5
$ws.Cells.Item(126, 16).Value = 1.5
$ws.Cells.Item(126, 17).Value = 1
$ws.Cells.Item(126, 18).Value = 2.025
$ws.Cells.Item(126, 19).Value = 1.775
$ws.Cells.Item(126, 20).Value = 3.25
$ws.Cells.Item(126, 21).Value = 1.95
$ws.Cells.Item(126, 22).Value = 1.85
$ws.Cells.Item(126, 23).Value = -1
$ws.Cells.Item(126, 24).Value = -1
$ws.Cells.Item(126, 25).Value = 0.5
$ws.Cells.Item(126, 26).Value = 0
$ws.Cells.Item(126, 27).Value = 0
$ws.Cells.Item(126, 28).Value = 0.95
$ws.Cells.Item(126, 29).Value = -1

# ---------------------------------------------------------------
# Row 127: Chennaiyin FC vs Northeast United - match result now known (2-1, Home win)
# ---------------------------------------------------------------
$ws.Cells.Item(127, 8).Value = 2
$ws.Cells.Item(127, 9).Value = 1
$ws.Cells.Item(127, 10).Value = 'H'
$ws.Cells.Item(127, 14).Value = 2
$ws.Cells.Item(127, 15).Value = 3.75
$ws.Cells.Item(127, 16).Value = 3.4
$ws.Cells.Item(127, 17).Value = -0.5
$ws.Cells.Item(127, 18).Value = 2.025
$ws.Cells.Item(127, 19).Value = 1.825
$ws.Cells.Item(127, 20).Value = 3
$ws.Cells.Item(127, 21).Value = 1.85
$ws.Cells.Item(127, 22).Value = 2
$ws.Cells.Item(127, 23).Value = 1
$ws.Cells.Item(127, 24).Value = -1
$ws.Cells.Item(127, 25).Value = -1
$ws.Cells.Item(127, 26).Value = 1.025
$ws.Cells.Item(127, 27).Value = -1
$ws.Cells.Item(127, 28).Value = 0
$ws.Cells.Item(127, 29).Value = 0

# ---------------------------------------------------------------
# Row 128: Punjab FC vs East Bengal Club - odds update (no result yet)
# ---------------------------------------------------------------
$ws.Cells.Item(128, 14).Value = 3.1
$ws.Cells.Item(128, 15).Value = 3.4
$ws.Cells.Item(128, 16).Value = 2.15
$ws.Cells.Item(128, 17).Value = 0.25
$ws.Cells.Item(128, 18).Value = 1.95
$ws.Cells.Item(128, 19).Value = 1.85
$ws.Cells.Item(128, 20).Value = 2.5
$ws.Cells.Item(128, 21).Value = 1.8
$ws.Cells.Item(128, 22).Value = 2

# ---------------------------------------------------------------
# Row 129 (new): Bengaluru vs Mohun Bagan SG
# ---------------------------------------------------------------
$ws.Range("A128").Copy()
$ws.Range("A129").PasteSpecial(-4122)
$ws.Range("E128").Copy()
$ws.Range("E129").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(129, 1).Value = 127
$ws.Cells.Item(129, 2).Value = 7749763
$ws.Cells.Item(129, 3).Value = 'India Super League'
$ws.Cells.Item(129, 4).Value = 'India Super League'
$ws.Cells.Item(129, 5).Value = 45393.45833333334
$ws.Cells.Item(129, 6).Value = 'Bengaluru'
$ws.Cells.Item(129, 7).Value = 'Mohun Bagan SG'
$ws.Cells.Item(129, 11).Value = 4.333
$ws.Cells.Item(129, 12).Value = 3.6
$ws.Cells.Item(129, 13).Value = 1.8
$ws.Cells.Item(129, 14).Value = 4.333
$ws.Cells.Item(129, 15).Value = 3.6
$ws.Cells.Item(129, 16).Value = 1.75
$ws.Cells.Item(129, 17).Value = 0.75
$ws.Cells.Item(129, 18).Value = 1.825
$ws.Cells.Item(129, 19).Value = 1.975
$ws.Cells.Item(129, 20).Value = 2.5
$ws.Cells.Item(129, 21).Value = 1.825
$ws.Cells.Item(129, 22).Value = 1.975
$ws.Cells.Item(129, 23).Value = 0
$ws.Cells.Item(129, 24).Value = 0
$ws.Cells.Item(129, 25).Value = 0
$ws.Cells.Item(129, 26).Value = 0
$ws.Cells.Item(129, 27).Value = 0

# ---------------------------------------------------------------
# Row 130 (new): Hyderabad FC vs Kerala Blasters
# ---------------------------------------------------------------
$ws.Range("A128").Copy()
$ws.Range("A130").PasteSpecial(-4122)
$ws.Range("E128").Copy()
$ws.Range("E130").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(130, 1).Value = 128
$ws.Cells.Item(130, 2).Value = 7749472
$ws.Cells.Item(130, 3).Value = 'India Super League'
$ws.Cells.Item(130, 4).Value = 'India Super League'
$ws.Cells.Item(130, 5).Value = 45394.45833333334
$ws.Cells.Item(130, 6).Value = 'Hyderabad FC'
$ws.Cells.Item(130, 7).Value = 'Kerala Blasters'
$ws.Cells.Item(130, 11).Value = 5.25
$ws.Cells.Item(130, 12).Value = 4
$ws.Cells.Item(130, 13).Value = 1.6
$ws.Cells.Item(130, 14).Value = 5.5
$ws.Cells.Item(130, 15).Value = 3.8
$ws.Cells.Item(130, 16).Value = 1.6
$ws.Cells.Item(130, 17).Value = 1
$ws.Cells.Item(130, 18).Value = 1.775
$ws.Cells.Item(130, 19).Value = 2.025
$ws.Cells.Item(130, 20).Value = 2.5
$ws.Cells.Item(130, 21).Value = 1.875
$ws.Cells.Item(130, 22).Value = 1.925
$ws.Cells.Item(130, 23).Value = 0
$ws.Cells.Item(130, 24).Value = 0
$ws.Cells.Item(130, 25).Value = 0
$ws.Cells.Item(130, 26).Value = 0
$ws.Cells.Item(130, 27).Value = 0
